$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price / volume / hour data as of the Feb 8, 2023 GitHub Actions run
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '331.88'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '1.04%'
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = '7'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '45.81'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '4.30%'
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '7'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.600'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.54%'
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '7'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08347'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '4.01%'
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = '7'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '3.02%'
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = '7'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9767'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '2.88%'
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = '7'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.572'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '0.30%'
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '7'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1163'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '2.78%'
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '7'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1927'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '3.55%'
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '7'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '10.39'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-3.10%'
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '7'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09926'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '0.12%'
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '7'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.04674'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.03%'
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '7'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.64%'
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '7'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001283'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.43%'
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '7'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006029'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.88%'
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '7'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.378'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.58%'
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '7'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.478'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '2.74%'
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '7'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-3.16%'
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '7'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1401'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-1.31%'
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '7'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.2652'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '4.09%'
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '7'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.04200'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '3.07%'
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '7'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.001312'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '3.92%'
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = '7'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.004650'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '7.36%'
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '7'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0001282'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '6.81%'
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = '7'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0003747'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '0.01%'
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '7'
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '7'
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '7'
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '7'
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '7'
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '7'
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '7'
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '7'
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = '7'
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '7'
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '7'
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '7'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02771'
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '7'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05793'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '2.69%'
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = '7'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.007842'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '3.73%'
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = '7'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1436'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '2.88%'
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = '7'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.007267'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-3.13%'
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = '7'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-0.12%'
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = '7'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008108'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-5.73%'
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = '7'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3405'
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = '7'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00007308'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '2.81%'
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = '7'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000751'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.13%'
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = '7'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '0.00%'
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = '7'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.003504'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-5.44%'
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = '7'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.003504'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.70%'
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = '7'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.13%'
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = '7'
